# Update "想去人数" (want-to-go count) figures in column F across the
# 展览 (exhibitions), 演出 (shows) and 全部类型 (all types) sheets, as
# output was regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5106
$ws1.Range("F4").Value = 7
$ws1.Range("F5").Value = 7397
$ws1.Range("F7").Value = 70
$ws1.Range("F12").Value = 4299
$ws1.Range("F13").Value = 1745
$ws1.Range("F16").Value = 2906
$ws1.Range("F18").Value = 565
$ws1.Range("F20").Value = 488
$ws1.Range("F21").Value = 429
$ws1.Range("F22").Value = 454
$ws1.Range("F23").Value = 300
$ws1.Range("F24").Value = 95
$ws1.Range("F26").Value = 1173
$ws1.Range("F28").Value = 1372
$ws1.Range("F29").Value = 106
$ws1.Range("F30").Value = 577
$ws1.Range("F32").Value = 513
$ws1.Range("F36").Value = 59
$ws1.Range("F37").Value = 2855
$ws1.Range("F38").Value = 699
$ws1.Range("F39").Value = 10
$ws1.Range("F40").Value = 48
$ws1.Range("F41").Value = 41
$ws1.Range("F42").Value = 13

# ---- Sheet "演出" (Shows) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 6
$ws2.Range("F3").Value = 8

# ---- Sheet "全部类型" (All types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5106
$ws4.Range("F4").Value = 7
$ws4.Range("F5").Value = 7397
$ws4.Range("F7").Value = 70
$ws4.Range("F12").Value = 4299
$ws4.Range("F13").Value = 1745
$ws4.Range("F16").Value = 2906
$ws4.Range("F18").Value = 565
$ws4.Range("F20").Value = 488
$ws4.Range("F21").Value = 429
$ws4.Range("F22").Value = 454
$ws4.Range("F23").Value = 6
$ws4.Range("F24").Value = 300
$ws4.Range("F25").Value = 95
$ws4.Range("F27").Value = 1173
$ws4.Range("F29").Value = 1372
$ws4.Range("F30").Value = 106
$ws4.Range("F31").Value = 577
$ws4.Range("F33").Value = 513
$ws4.Range("F37").Value = 59
$ws4.Range("F38").Value = 2855
$ws4.Range("F39").Value = 8
$ws4.Range("F40").Value = 699
$ws4.Range("F41").Value = 10
$ws4.Range("F42").Value = 48
$ws4.Range("F43").Value = 41
$ws4.Range("F44").Value = 13
